$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1498.3334
$ws.Range("I18").Value = 1497.5
$ws.Range("K18").Value = 1497.5
$ws.Range("M18").Value = -1213.5
$ws.Range("H33").Value = 139.73334
$ws.Range("I33").Value = 127.72727
$ws.Range("K33").Value = 127.72727
$ws.Range("M33").Value = 101.27273
$ws.Range("H53").Value = 1840.2858
$ws.Range("I53").Value = 2303.2
$ws.Range("J53").Value = 683
$ws.Range("K53").Value = 2303.2
$ws.Range("L53").Value = 683
$ws.Range("M53").Value = -1666.2
$ws.Range("N53").Value = -1957
$ws.Range("H80").Value = 1086.5385
$ws.Range("J80").Value = 1041.6666
$ws.Range("L80").Value = 3124.9998
$ws.Range("N80").Value = -5120.9998
$ws.Range("H83").Value = 1086.5385
$ws.Range("J83").Value = 1041.6666
$ws.Range("L83").Value = 9374.999400000001
$ws.Range("N83").Value = -19358.9994
$ws.Range("H87").Value = 56666.668
$ws.Range("J87").Value = 56666.668
$ws.Range("L87").Value = 56666.668
$ws.Range("N87").Value = -59162.668
$ws.Range("H90").Value = 56666.668
$ws.Range("J90").Value = 56666.668
$ws.Range("L90").Value = 170000.004
$ws.Range("N90").Value = -182480.004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1380.75
$ws.Range("I2").Value = 1380.75
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1380.75
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1267.75
$ws.Range("N2").ClearContents()
$ws.Range("H5").Value = 246.75
$ws.Range("I5").Value = 128.33333
$ws.Range("K5").Value = 128.33333
$ws.Range("M5").Value = -16.33332999999999
$ws.Range("H74").Value = 1735.5454
$ws.Range("I74").Value = 1699.375
$ws.Range("J74").Value = 1832
$ws.Range("K74").Value = 1699.375
$ws.Range("L74").Value = 1832
$ws.Range("M74").Value = -825.375
$ws.Range("N74").Value = -3580
$ws.Range("H77").Value = 1735.5454
$ws.Range("I77").Value = 1699.375
$ws.Range("J77").Value = 1832
$ws.Range("K77").Value = 8496.875
$ws.Range("L77").Value = 9160
$ws.Range("M77").Value = -4128.875
$ws.Range("N77").Value = -17896
$ws.Range("H116").Value = 1380.75
$ws.Range("I116").Value = 1380.75
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1380.75
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 913.25
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1380.75
$ws.Range("I3").Value = 1380.75
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1380.75
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1266.75
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 246.75
$ws.Range("I4").Value = 128.33333
$ws.Range("K4").Value = 128.33333
$ws.Range("M4").Value = -13.33332999999999
$ws.Range("H20").Value = 1365.2858
$ws.Range("I20").Value = 1441.6
$ws.Range("J20").Value = 1174.5
$ws.Range("K20").Value = 1441.6
$ws.Range("L20").Value = 1174.5
$ws.Range("M20").Value = -1194.6
$ws.Range("N20").Value = -1668.5
$ws.Range("H82").Value = 21327.46
$ws.Range("H85").Value = 21327.46
$ws.Range("H86").Value = 3100.2
$ws.Range("I86").Value = 3375.375
$ws.Range("J86").Value = 2785.7144
$ws.Range("K86").Value = 3375.375
$ws.Range("L86").Value = 2785.7144
$ws.Range("M86").Value = -2252.375
$ws.Range("N86").Value = -5031.7144
$ws.Range("H89").Value = 3100.2
$ws.Range("I89").Value = 3375.375
$ws.Range("J89").Value = 2785.7144
$ws.Range("K89").Value = 16876.875
$ws.Range("L89").Value = 13928.572
$ws.Range("M89").Value = -11260.875
$ws.Range("N89").Value = -25160.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 650.8182
$ws.Range("I105").Value = 650.8182
$ws.Range("K105").Value = 650.8182
$ws.Range("M105").Value = 1096.1818
$ws.Range("H134").Value = 2733.4285
$ws.Range("I134").Value = 2599.923
$ws.Range("J134").Value = 4469
$ws.Range("K134").Value = 7799.768999999999
$ws.Range("L134").Value = 13407
$ws.Range("M134").Value = -5264.768999999999
$ws.Range("N134").Value = -18477

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 7423.6665
$ws.Range("I81").Value = 7388
$ws.Range("K81").Value = 22164
$ws.Range("M81").Value = -21041
$ws.Range("H84").Value = 7423.6665
$ws.Range("I84").Value = 7388
$ws.Range("K84").Value = 66492
$ws.Range("M84").Value = -60876

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1505000
$ws.Range("I70").Value = 1505000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 1505000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -1504730
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 1505000
$ws.Range("I73").Value = 1505000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 1505000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1504064
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 2483
$ws.Range("J80").Value = 2379.8
$ws.Range("L80").Value = 2379.8
$ws.Range("N80").Value = -4375.8
$ws.Range("H83").Value = 2483
$ws.Range("J83").Value = 2379.8
$ws.Range("L83").Value = 11899
$ws.Range("N83").Value = -21883
$ws.Range("H126").Value = 6216.3335
$ws.Range("I126").Value = 7199.6
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 21598.8
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -19128.8
$ws.Range("N126").Value = -8840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 266.66666
$ws.Range("I55").Value = 250.6
$ws.Range("J55").Value = 347
$ws.Range("K55").Value = 250.6
$ws.Range("L55").Value = 347
$ws.Range("M55").Value = -77.59999999999999
$ws.Range("N55").Value = -693
$ws.Range("H61").Value = 7833.3335
$ws.Range("I61").Value = 7833.3335
$ws.Range("K61").Value = 7833.3335
$ws.Range("M61").Value = -7631.3335
$ws.Range("H82").Value = 2006
$ws.Range("I82").Value = 1924.8334
$ws.Range("J82").Value = 2249.5
$ws.Range("K82").Value = 1924.8334
$ws.Range("L82").Value = 2249.5
$ws.Range("M82").Value = -1563.8334
$ws.Range("N82").Value = -2971.5
$ws.Range("H85").Value = 2006
$ws.Range("I85").Value = 1924.8334
$ws.Range("J85").Value = 2249.5
$ws.Range("K85").Value = 1924.8334
$ws.Range("L85").Value = 2249.5
$ws.Range("M85").Value = -676.8334
$ws.Range("N85").Value = -4745.5
$ws.Range("H113").Value = 7833.3335
$ws.Range("I113").Value = 7833.3335
$ws.Range("K113").Value = 7833.3335
$ws.Range("M113").Value = -5663.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 9250.5
$ws.Range("I13").Value = 9250.5
$ws.Range("K13").Value = 9250.5
$ws.Range("M13").Value = -9110.5
$ws.Range("H109").Value = 75000
$ws.Range("J109").Value = 75000
$ws.Range("L109").Value = 75000
$ws.Range("N109").Value = -77774
$ws.Range("H113").Value = 373.6
$ws.Range("I113").Value = 359.55554
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 1078.66662
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1091.33338
$ws.Range("N113").Value = -5840
$ws.Range("H126").Value = 22275.133
$ws.Range("I126").Value = 23273.572
$ws.Range("J126").Value = 8297
$ws.Range("K126").Value = 69820.716
$ws.Range("L126").Value = 24891
$ws.Range("M126").Value = -67350.716
$ws.Range("N126").Value = -29831
